{"js": "// Applies the textual corrections described in the commit\n// \"Entrega relatorio de Proposta\" to the document body.\n// Each entry is a unique (search, replacement) pair; the search text\n// is located with Body.search (case-sensitive, punctuation-matching)\n// and the matched range's text is swapped in place via insertText with\n// Word.InsertLocation.replace, which preserves the surrounding run\n// formatting.\n\nconst edits = [\n  {\n    find: \"a entidade representante da ind\u00fastria, via contacto movel, correio eletr\u00f3nico ou\",\n    replace: \"a entidade representante da Ind\u00fastria, via contacto m\u00f3vel, correio eletr\u00f3nico ou\"\n  },\n  {\n    find: \"Esse interesse traduz-se no objetivo de fornecer aos seus clientes e empresas um estoque de produtos\",\n    replace: \"Esse interesse traduz-se no objetivo de fornecer aos seus clientes e empresas um stock de produtos\"\n  },\n  {\n    find: \"Quando um cliente entra em contato com a ind\u00fastria fornecedora,\",\n    replace: \"Quando um cliente entra em contacto com a Ind\u00fastria fornecedora,\"\n  },\n  {\n    find: \"Os trabalhadores da ind\u00fastria se dedicam \u00e0s suas tarefas\",\n    replace: \"Os trabalhadores da ind\u00fastria dedicam-se \u00e0s suas tarefas\"\n  },\n  {\n    find: \"para equipas de ge\",\n    replace: \"para as equipas de ge\"\n  },\n  {\n    find: \"todos riscos e benef\u00edcios\",\n    replace: \"todos os riscos e benef\u00edcios\"\n  },\n  {\n    find: \"adicionado \u00e1 base de dados\",\n    replace: \"adicionado \u00e0 base de dados\"\n  },\n  {\n    find: \"na pr\u00f3pria loja ou solicitar a entrega\",\n    replace: \"na pr\u00f3pria fornecedora ou solicitar a entrega\"\n  },\n  {\n    find: \"ser um levantamento na fabrica,\",\n    replace: \"ser um levantamento na f\u00e1brica,\"\n  },\n  {\n    find: \"O representante da ind\u00fastria possui acesso ao fluxo de encomendas,\",\n    replace: \"O representante da Ind\u00fastria possui acesso ao fluxo de encomendas,\"\n  }\n];\n\nfor (const edit of edits) {\n  const results = context.document.body.search(edit.find, { matchCase: true });\n  results.load(\"text\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + edit.find);\n  }\n\n  results.items[0].insertText(edit.replace, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Applies the textual corrections described in the commit\n# \"Entrega relatorio de Proposta\" using Word's Find & Replace\n# (Range.Find.Execute), which preserves the surrounding run formatting\n# and merges the replaced text into adjacent identically-formatted runs,\n# matching how Word itself performs a find/replace edit.\n\n$d = $word.ActiveDocument\n\nfunction Replace-Once($findText, $replaceText) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2) | Out-Null\n}\n\nReplace-Once \"a entidade representante da ind\u00fastria, via contacto movel, correio eletr\u00f3nico ou\" \"a entidade representante da Ind\u00fastria, via contacto m\u00f3vel, correio eletr\u00f3nico ou\"\nReplace-Once \"Esse interesse traduz-se no objetivo de fornecer aos seus clientes e empresas um estoque de produtos\" \"Esse interesse traduz-se no objetivo de fornecer aos seus clientes e empresas um stock de produtos\"\nReplace-Once \"Quando um cliente entra em contato com a ind\u00fastria fornecedora,\" \"Quando um cliente entra em contacto com a Ind\u00fastria fornecedora,\"\nReplace-Once \"Os trabalhadores da ind\u00fastria se dedicam \u00e0s suas tarefas\" \"Os trabalhadores da ind\u00fastria dedicam-se \u00e0s suas tarefas\"\nReplace-Once \"para equipas de ge\" \"para as equipas de ge\"\nReplace-Once \"todos riscos e benef\u00edcios\" \"todos os riscos e benef\u00edcios\"\nReplace-Once \"adicionado \u00e1 base de dados\" \"adicionado \u00e0 base de dados\"\nReplace-Once \"na pr\u00f3pria loja ou solicitar a entrega\" \"na pr\u00f3pria fornecedora ou solicitar a entrega\"\nReplace-Once \"ser um levantamento na fabrica,\" \"ser um levantamento na f\u00e1brica,\"\nReplace-Once \"O representante da ind\u00fastria possui acesso ao fluxo de encomendas,\" \"O representante da Ind\u00fastria possui acesso ao fluxo de encomendas,\"\n"}
